$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$down = [char]8600
$up   = [char]8599

$data = @{
    2 = @("BINANCE_SPOT_ADA_USDT", "0/0    $down", "0/0    $down", "0/0    $down", "0/0    $up",   "0/0    $down")
    3 = @("BINANCE_SPOT_BTC_USDT", "1/0    $down", "5/1    $down", "6/1    $down", "1/0    $up",   "0/0    $down")
    4 = @("BINANCE_SPOT_ETH_USDT", "1/0    $down", "0/0    $down", "5/1    $down", "0/0    $down", "0/0    $down")
    5 = @("BINANCE_SPOT_ADA_USDT", "0/0    $down", "0/0    $down", "0/0    $down", "0/0    $up",   "0/0    $down")
    6 = @("BINANCE_SPOT_SOL_USDT", "9/2    $up",   "0/0    $up",   "0/0    $up",   "0/0    $up",   "0/0    $up")
    7 = @("BINANCE_SPOT_CRV_USDT", "0/0    $down", "0/0    $up",   "0/0    $down", "0/0    $up",   "2/0    $up")
    8 = @("BINANCE_SPOT_BTC_USDT", "1/0    $down", "5/1    $down", "6/1    $down", "1/0    $up",   "0/0    $down")
}

foreach ($rowNum in $data.Keys) {
    $values = $data[$rowNum]
    $ws.Cells.Item($rowNum, 1).Value = $values[0]
    $ws.Cells.Item($rowNum, 2).Value = $values[1]
    $ws.Cells.Item($rowNum, 3).Value = $values[2]
    $ws.Cells.Item($rowNum, 4).Value = $values[3]
    $ws.Cells.Item($rowNum, 5).Value = $values[4]
    $ws.Cells.Item($rowNum, 6).Value = $values[5]
}
